# Marksheet update: fill in actual exam results (was template/"Absent" placeholder).
# "Handles float input without breaking stuff" - scores now reflect a real attempt
# instead of the all-zero / Absent template, and the stray 3rd Student/Correct Ans
# block (columns G:H) is removed since only two answer-pairs are used.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary rows (10-12): give the title cells the same "mtitleStyle" as other
#     row labels, and fill in the real Right/Wrong/NotAttempt/Max + Marking + Total figures.
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("A12").Style = "mtitleStyle"

$ws.Range("B10").Value = 19
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 6
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1
$ws.Range("D11").Value = 0

$ws.Range("B12").Value = 76
$ws.Range("C12").Value = -3
$ws.Range("E12").Value = "73/112"

# --- Per-question "Student Ans" columns (A and D). Correct answers (B/E) are
#     already present; we only need to mark what the student actually answered,
#     using correctStyle/incorrectStyle/normalStyle (blank = not attempted).

$ws.Range("A16").Value = "Option A"
$ws.Range("A16").Style = "correctStyle"

$ws.Range("D16").Value = "Option A"
$ws.Range("D16").Style = "correctStyle"

$ws.Range("D17").Value = "Option C"
$ws.Range("D17").Style = "correctStyle"

$ws.Range("A18").Value = "Option B"
$ws.Range("A18").Style = "correctStyle"

$ws.Range("D18").Value = "Option D"
$ws.Range("D18").Style = "correctStyle"

$ws.Range("A19").Value = "Option C"
$ws.Range("A19").Style = "correctStyle"

$ws.Range("A20").Value = "Option B"
$ws.Range("A20").Style = "correctStyle"

$ws.Range("A21").Value = "Option C"
$ws.Range("A21").Style = "correctStyle"

$ws.Range("A22").Value = "Option D"
$ws.Range("A22").Style = "correctStyle"

$ws.Range("A26").Value = "Option D"
$ws.Range("A26").Style = "incorrectStyle"

$ws.Range("A27").Value = "Option A"
$ws.Range("A27").Style = "correctStyle"

$ws.Range("A28").Value = "Option D"
$ws.Range("A28").Style = "correctStyle"

$ws.Range("A30").Value = "Option B"
$ws.Range("A30").Style = "correctStyle"

$ws.Range("A31").Value = "Option B"
$ws.Range("A31").Style = "incorrectStyle"

$ws.Range("A32").Value = "Option C"
$ws.Range("A32").Style = "correctStyle"

$ws.Range("A33").Value = "Option D"
$ws.Range("A33").Style = "correctStyle"

$ws.Range("A35").Value = "Option D"
$ws.Range("A35").Style = "correctStyle"

$ws.Range("A36").Value = "Option D"
$ws.Range("A36").Style = "incorrectStyle"

$ws.Range("A37").Value = "Option A"
$ws.Range("A37").Style = "correctStyle"

$ws.Range("A38").Value = "Option A"
$ws.Range("A38").Style = "correctStyle"

$ws.Range("A39").Value = "Option D"
$ws.Range("A39").Style = "correctStyle"

$ws.Range("A40").Value = "Option D"
$ws.Range("A40").Style = "correctStyle"

# --- Only the first 3 questions (rows 16-18) actually use a second Student/Correct
#     Ans pair (columns D:E); rows 19-40 never did, so wipe those leftover D:E cells.
$ws.Range("D19:E40").Clear()

# --- Drop the unused 3rd Student/Correct-Ans block (columns G:H) entirely so the
#     sheet's used range shrinks from A5:H40 down to A5:E40.
$ws.Range("G1:H1").EntireColumn.Delete()
